$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10; this shifts rows 10..101 down to 11..102
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with data
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(10, 3).Value = "Coquimbo"
$ws.Cells.Item(10, 4).Value = 45083
$ws.Cells.Item(10, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(10, 5).Value = 4
$ws.Cells.Item(10, 6).Value = "Fruta"
$ws.Cells.Item(10, 7).Value = 100104
$ws.Cells.Item(10, 8).Value = "Frutos de pepita"
$ws.Cells.Item(10, 9).Value = 100104003
$ws.Cells.Item(10, 10).Value = "Membrillo"
$ws.Cells.Item(10, 11).Value = "Champion"
$ws.Cells.Item(10, 12).Value = "Primera"
$ws.Cells.Item(10, 13).Value = 14
$ws.Cells.Item(10, 14).Value = 200000
$ws.Cells.Item(10, 15).Value = 210000
$ws.Cells.Item(10, 16).Value = 205000
$ws.Cells.Item(10, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(10, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(10, 19).Value = 456
$ws.Cells.Item(10, 20).Value = 450
